$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates - values that are NOT a single valid number
# (e.g. contain multiple dots like "26.253.16") are plain text already.
$ws.Range("D2").Value = "26.253.16"
$ws.Range("D3").Value = "1.592.24"
$ws.Range("D12").Value = "1.816.48"
$ws.Range("D13").Value = "1.612.08"
$ws.Range("D17").Value = "26.255.42"
$ws.Range("D33").Value = "1.418.76"
$ws.Range("D45").Value = "1.728.67"

# Price column (D) updates - these values parse as a plain decimal number,
# so a leading apostrophe is used to force Excel to keep them as text,
# matching the sheet's existing inline-string price formatting.
$ws.Range("D5").Value = "'212.96"
$ws.Range("D10").Value = "'18.97"
$ws.Range("D15").Value = "'0.508"
$ws.Range("D16").Value = "'63.82"
$ws.Range("D19").Value = "'215.69"
$ws.Range("D20").Value = "'7.36"
$ws.Range("D22").Value = "'4.29"
$ws.Range("D25").Value = "'145.28"
$ws.Range("D30").Value = "'0.0494"
$ws.Range("D37").Value = "'0.573"
$ws.Range("D39").Value = "'0.825"
$ws.Range("D40").Value = "'5.78"
$ws.Range("D42").Value = "'0.938"
$ws.Range("D44").Value = "'0.761"
$ws.Range("D46").Value = "'60.90"
$ws.Range("D47").Value = "'86.88"
$ws.Range("D48").Value = "'1.48"
$ws.Range("D50").Value = "'0.0953"

# Volume(1h) column (E) updates
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("E15").Value = "  -2.33%  "
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("E30").Value = "  -1.33%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("E33").Value = "  +5.90%  "
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  -9.92%  "
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("E50").Value = "  -2.57%  "
$ws.Range("E51").Value = "  +0.10%  "
